$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 68 ---
$ws.Range("A68").Value = 45643
$ws.Range("B68").Value = 597
$ws.Range("D68").Value = 1
$ws.Range("E68").Value = "traded using OI data and got good result"

# --- Row 69 ---
$ws.Range("A69").Value = 45644
$ws.Range("B69").Value = 600
$ws.Range("D69").Value = 1
$ws.Range("E69").Value = "traded using OI data and got good result"

# --- Row 70 ---
$ws.Range("A70").Value = 45645
$ws.Range("B70").Value = 370
$ws.Range("D70").Value = 1
$ws.Range("E70").Value = "traded using OI data but not sure about market trend so exited early"

# Copy the date number format from the previous row (A67) onto the new date cells
$ws.Range("A67").Copy()
$ws.Range("A68:A70").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the visible window / selection to mirror the saved view state
$ws.Application.ActiveWindow.ScrollRow = 60
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B70").Select()
